$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# --- Edit 1: "Text Placeholder 3" -------------------------------------
# Merge the three runs ("Insight ", "into most ", "valuable suburbs to
# invest in") into a single run with the combined text.
$titleShape = $s.Shapes.Item("Text Placeholder 3")
$titleTr = $titleShape.TextFrame.TextRange
# First set to a distinct placeholder so the subsequent assignment of the
# (coincidentally identical-looking) concatenated text is not treated as
# a no-op by the engine.
$titleTr.Text = "TEMP_PLACEHOLDER_TEXT"
$titleShape.TextFrame.TextRange.Text = "Insight into most valuable suburbs to invest in"

# --- Edit 2: "Content Placeholder 5" ----------------------------------
# Add two new paragraphs of analysis text before the existing (blank)
# trailing paragraph.
$contentShape = $s.Shapes.Item("Content Placeholder 5")
$contentTr = $contentShape.TextFrame.TextRange

# Insert the first new paragraph's text right before the existing blank
# paragraph (no leading/trailing CR yet -> merges with the existing
# endParaRPr for now).
$contentTr.InsertBefore("Using the choropleth map, it is clear that many Melbourne suburbs are experiencing growth at the moment")

# Re-assign the whole range via a temp value first (same no-op guard as
# above), then to the real text, so the run picks up a proper <a:rPr>.
$wholeTr = $contentShape.TextFrame.TextRange
$wholeTr.Text = "TEMP_PLACEHOLDER_TEXT_1"
$contentShape.TextFrame.TextRange.Text = "Using the choropleth map, it is clear that many Melbourne suburbs are experiencing growth at the moment"

# Append a second paragraph (temp text) followed by a paragraph break so
# the original blank paragraph/endParaRPr is preserved as a third,
# trailing, empty paragraph.
$afterTr = $contentShape.TextFrame.TextRange
$afterTr.InsertAfter("`rTEMP_PLACEHOLDER_TEXT_2`r")

# Replace the second paragraph's temp text with the real text.
$para2 = $contentShape.TextFrame.TextRange.Paragraphs(2, 1)
$para2.Text = "Many of the areas with the highest growth rates appear to be in the northern and eastern suburbs, signaling them as great opportunities for investment based on past performance"

# The trailing third paragraph currently holds a stray empty run created
# by the paragraph-break insertion above; clear its text so it collapses
# back down to just the original <a:endParaRPr/>.
$para3 = $contentShape.TextFrame.TextRange.Paragraphs(3, 1)
$run3 = $para3.Runs(1, 1)
$run3.Text = ""
